$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24:46 down to 25:47
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with its data
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44790
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112040
$ws.Range("G24").Value = "Cilantro"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 10000
$ws.Range("M24").Value = 10000
$ws.Range("N24").Value = '$/caja 36 atados'
$ws.Range("O24").Value = "Provincia de Quillota"
$ws.Range("P24").Value = 278
$ws.Range("Q24").Value = 36
$ws.Range("R24").Value = "Hortaliza"
